# Auto-generated Excel COM-interop edit script
# Applies the diff: updates "today" (row 13) and MTD (row 14) KPI rows
# on Dashboard, Kerala, Tamilnadu and Chennai sheets, and refreshes / extends
# the per-hour breakdown on the "Hourly Report" sheet (rows 280-283 updated,
# rows 284-288 added for hours 19-23).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Dashboard / Kerala / Tamilnadu / Chennai -- row 13 ("today") + row 14 (MTD)
# ---------------------------------------------------------------------------
$sheetUpdates = @(
  @{
    Name = "Dashboard"
    Row13 = @{ "B"=21240; "C"=11411; "D"=2583; "E"=8828; "F"=0.22636052931382; "G"=45.26001226886338; "H"=9829; "I"=0.537241054613936; "J"=9829; "K"=9755; "L"=9318; "M"=9148; "N"=170; "O"=0.9552024602767811; "P"=437; "Q"=0.04479753972321888; "R"=74; "S"=8734; "T"=101; "U"=584; "V"=0.9272746576069646; "W"=133.0100568430258; "X"=18.45376038478356; "Y"=693; "Z"=58; "AA"=45; "AE"=670; "AF"=58; "AG"=44; "AH"=0.9668109668109668; "AI"=1; "AJ"=0.9777777777777777 }
    Row14 = @{ "B"=234697; "C"=122037; "D"=30579; "E"=91458; "F"=0.2505715479731557; "G"=43.05134508386801; "H"=112660; "I"=0.5199768211779443; "J"=112660; "K"=111747; "L"=105383; "M"=103355; "N"=2028; "O"=0.943049925277636; "P"=6364; "Q"=0.05695007472236398; "R"=913; "S"=100412; "T"=1989; "U"=4971; "V"=0.9351786313005251; "W"=133.3651492429007; "X"=17.17855933433312; "Y"=6952; "Z"=386; "AA"=371; "AE"=6814; "AF"=386; "AG"=369; "AH"=0.9801495972382048; "AJ"=0.9946091644204852 }
  },
  @{
    Name = "Kerala"
    Row13 = @{ "B"=13784; "C"=8108; "D"=1485; "E"=6623; "F"=0.1831524420325604; "G"=52.15355204736063; "H"=5676; "I"=0.5882182240278584; "J"=5676; "K"=5631; "L"=5312; "M"=5240; "N"=72; "O"=0.9433493162848517; "P"=319; "Q"=0.0566506837151482; "R"=45; "S"=4775; "T"=82; "U"=537; "V"=0.8852428624397478; "W"=133.5335877862595; "X"=26.24541984732824; "Y"=450; "Z"=50; "AA"=31; "AE"=427; "AF"=50; "AG"=30; "AH"=0.9488888888888889; "AI"=1; "AJ"=0.967741935483871 }
    Row14 = @{ "B"=134947; "C"=80514; "D"=17164; "E"=63350; "F"=0.2131803164667014; "G"=48.83791638721216; "H"=54433; "I"=0.5966342341808265; "J"=54433; "K"=54076; "L"=51515; "M"=50661; "N"=854; "O"=0.9526407278644871; "P"=2561; "Q"=0.04735927213551292; "R"=357; "S"=48070; "T"=648; "U"=3445; "V"=0.9215344209497153; "W"=136.0433272142279; "X"=20.95325793016324; "Y"=3203; "Z"=325; "AA"=265; "AE"=3128; "AF"=325; "AG"=263; "AH"=0.9765844520761786; "AJ"=0.9924528301886792 }
  },
  @{
    Name = "Tamilnadu"
    Row13 = @{ "B"=4566; "C"=2126; "D"=745; "E"=1381; "F"=0.3504233301975541; "G"=25.91768579492004; "H"=2440; "I"=0.4656154183092422; "J"=2440; "K"=2415; "L"=2342; "M"=2281; "N"=61; "O"=0.9697722567287783; "P"=73; "Q"=0.03022774327122164; "R"=25; "S"=2315; "T"=11; "U"=27; "V"=0.9838504037399065; "W"=128.8728627794827; "X"=7.828583954405962; "Y"=154; "Z"=1; "AA"=11; "AE"=154; "AF"=1; "AG"=11; "AH"=1; "AI"=1; "AJ"=1 }
    Row14 = @{ "B"=60407; "C"=26816; "D"=9149; "E"=17667; "F"=0.3411769093078759; "G"=29.10012678997613; "H"=33591; "I"=0.4439220620126807; "J"=33591; "K"=33220; "L"=30860; "M"=30073; "N"=787; "O"=0.9289584587597832; "P"=2360; "Q"=0.07104154124021676; "R"=371; "S"=29983; "T"=805; "U"=877; "V"=0.9468814148113058; "W"=127.7579223888538; "X"=13.64020882519203; "Y"=2246; "Z"=21; "AA"=80; "AE"=2206; "AF"=21; "AG"=80; "AH"=0.9821905609973285 }
  },
  @{
    Name = "Chennai"
    Row13 = @{ "B"=2890; "C"=1177; "D"=353; "E"=824; "F"=0.2999150382327953; "G"=32.71028037383178; "H"=1713; "I"=0.4072664359861592; "J"=1713; "K"=1709; "L"=1664; "M"=1627; "N"=37; "O"=0.97366881217086; "P"=45; "Q"=0.02633118782913996; "S"=1644; "T"=8; "U"=20; "V"=0.9832535885167464; "W"=137.1241548862938; "X"=8.255685310387216; "Y"=89; "Z"=7; "AA"=3; "AE"=89; "AF"=7; "AG"=3; "AH"=1; "AI"=1; "AJ"=1 }
    Row14 = @{ "B"=39343; "C"=14707; "D"=4266; "E"=10441; "F"=0.2900659549874209; "G"=36.81049840212144; "H"=24636; "I"=0.3738149098950258; "J"=24636; "K"=24451; "L"=23008; "M"=22621; "N"=387; "O"=0.9409840088339945; "P"=1443; "Q"=0.0590159911660055; "S"=22359; "T"=536; "U"=649; "V"=0.9496687054026504; "W"=134.8216259228151; "X"=13.42889350603422; "Y"=1503; "Z"=40; "AA"=26; "AE"=1480; "AF"=40; "AG"=26; "AH"=0.9846972721224219 }
  },
)

foreach ($entry in $sheetUpdates) {
  $ws = $wb.Worksheets.Item($entry.Name)
  foreach ($col in $entry.Row13.Keys) {
    $ws.Range("$($col)13").Value = $entry.Row13[$col]
  }
  foreach ($col in $entry.Row14.Keys) {
    $ws.Range("$($col)14").Value = $entry.Row14[$col]
  }
}

# ---------------------------------------------------------------------------
# 2) Hourly Report -- update existing hours 15-18 (rows 280-283), then append
#    hours 19-23 (new rows 284-288), copying the formatting of row 283 down.
# ---------------------------------------------------------------------------
$hourly = $wb.Worksheets.Item("Hourly Report")

$hourlyUpdates = @{
  280 = @{ "C"=1085 }
  281 = @{ "C"=1177; "D"=515; "E"=508; "G"=0.9864077669902913; "H"=0.9941060903732809; "I"=153.3400809716599; "J"=99.6551724137931 }
  282 = @{ "C"=1378; "D"=624; "E"=619; "G"=0.9919871794871795; "H"=0.9903381642512076; "I"=141.9457236842105 }
  283 = @{ "C"=1569; "D"=709; "E"=697; "F"=12; "G"=0.9830747531734838; "H"=0.9497847919655668; "I"=146.3450292397661; "J"=91.48418491484185; "K"=100; "L"=100 }
  284 = @{ "A"=45912; "B"=19; "C"=1824; "D"=800; "E"=796; "F"=4; "G"=0.995; "H"=0.9861809045226131; "I"=118.9898089171975; "J"=97.92843691148776; "K"=100; "L"=100 }
  285 = @{ "A"=45912; "B"=20; "C"=1796; "D"=811; "E"=732; "F"=79; "G"=0.9025893958076449; "H"=0.8466666666666667; "I"=116.8022130013831; "J"=77.37373737373737; "K"=98.65771812080537; "L"=99.05660377358491 }
  286 = @{ "A"=45912; "B"=21; "C"=1524; "D"=682; "E"=587; "F"=95; "G"=0.8607038123167157; "H"=0.6891447368421053; "I"=134.5989583333333; "J"=56.90866510538641; "K"=97.34513274336283; "L"=97.05882352941177 }
  287 = @{ "A"=45912; "B"=22; "C"=878; "D"=407; "E"=380; "F"=27; "G"=0.9336609336609336; "H"=0.8756476683937824; "I"=139.5067385444744; "J"=89.24731182795699; "K"=88.60759493670885; "L"=67.85714285714286 }
  288 = @{ "A"=45912; "B"=23; "C"=435; "D"=213; "E"=213; "F"=0; "G"=1; "H"=1; "I"=131.8146341463415; "J"=100; "K"=100; "L"=100 }
}

# Rows 284-288 do not exist yet: clone formatting from row 283 first so the
# new rows inherit the same styles (date / number / percent) as the rest of
# the table, matching the workbook s="2..5" style pattern.
$hourly.Range("A283:L283").Copy($hourly.Range("A284:L284"))
$hourly.Range("A283:L283").Copy($hourly.Range("A285:L285"))
$hourly.Range("A283:L283").Copy($hourly.Range("A286:L286"))
$hourly.Range("A283:L283").Copy($hourly.Range("A287:L287"))
$hourly.Range("A283:L283").Copy($hourly.Range("A288:L288"))

foreach ($rowNum in $hourlyUpdates.Keys) {
  $rowData = $hourlyUpdates[$rowNum]
  foreach ($col in $rowData.Keys) {
    $hourly.Range("$($col)$($rowNum)").Value = $rowData[$col]
  }
}

